$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1190320826869504
$ws.Range("C2").Value = 10.34677158129881
$ws.Range("D2").Value = 22.3905356188092
$ws.Range("E2").Value = 91228006295.30009
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 91228006328.15643
